$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "92.245.98"
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").Value = "3.119.91"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "'243.52"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'617.23"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'1.10"
$ws.Range("E7").Value = "  -3.65%  "
$ws.Range("D8").Value = "'0.400"
$ws.Range("E8").Value = "  +10.51%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.118.83"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").Value = "'0.736"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "'0.203"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("E13").Value = "  +5.27%  "
$ws.Range("D14").Value = "'34.78"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "92.282.86"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("B16").Value = "Toncoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D16").Value = "'5.53"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").Value = "3.693.83"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "3.134.00"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "'3.69"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").Value = "'14.86"
$ws.Range("E20").Value = "  +3.49%  "
$ws.Range("D21").Value = "'5.81"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").Value = "'0.0000207"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").Value = "'450.11"
$ws.Range("E23").Value = "  +3.28%  "
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "'9.32"
$ws.Range("E24").Value = "  +4.22%  "
$ws.Range("D25").Value = "'5.66"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "'11.72"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'81.59"
$ws.Range("E27").Value = "  -9.88%  "
$ws.Range("D28").Value = "3.278.93"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").Value = "'0.138"
$ws.Range("E30").Value = "  +19.11%  "
$ws.Range("D31").Value = "'0.228"
$ws.Range("E31").Value = "  -5.68%  "
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("D33").Value = "'9.37"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("D34").Value = "'0.173"
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'8.02"
$ws.Range("E36").Value = "  +5.75%  "
$ws.Range("D37").Value = "'26.45"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Value = "'4.10"
$ws.Range("E38").Value = "  -4.16%  "
$ws.Range("D39").Value = "'1.94"
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("D40").Value = "'490.74"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").Value = "'1.31"
$ws.Range("E41").Value = "  +3.07%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'3.51"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.438"
$ws.Range("E43").Value = "  +5.78%  "
$ws.Range("D44").Value = "'22.22"
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'159.13"
$ws.Range("E46").Value = "  +3.19%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.94"
$ws.Range("E47").Value = "  +3.57%  "
$ws.Range("D48").Value = "'0.700"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("D49").Value = "'1.37"
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("D50").Value = "'0.0331"
$ws.Range("E50").Value = "  +7.58%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "'44.12"
$ws.Range("E51").Value = "  +0.10%  "
